$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "US Core PMO ServiceRequest Profile" row (row 43), shifting
# all subsequent rows up by one.
$ws.Rows.Item(43).Delete()
